$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "44.142.06"
Set-TextCell "E2" "  +0.98%  "

# Row 3
Set-TextCell "D3" "2.248.64"
Set-TextCell "E3" "  +0.09%  "

# Row 4
Set-TextCell "E4" "  +0.23%  "

# Row 5
Set-TextCell "D5" "306.46"
Set-TextCell "E5" "  -4.97%  "

# Row 6
Set-TextCell "D6" "97.79"
Set-TextCell "E6" "  -3.35%  "

# Row 7
Set-TextCell "D7" "0.575"
Set-TextCell "E7" "  -0.78%  "

# Row 8
Set-TextCell "E8" "  +0.15%  "

# Row 9
Set-TextCell "D9" "0.531"
Set-TextCell "E9" "  -4.22%  "

# Row 10
Set-TextCell "D10" "35.33"
Set-TextCell "E10" "  -4.23%  "

# Row 11
Set-TextCell "D11" "0.0817"
Set-TextCell "E11" "  -1.48%  "

# Row 12
Set-TextCell "D12" "7.26"
Set-TextCell "E12" "  -6.12%  "

# Row 13
Set-TextCell "E13" "  -2.00%  "

# Row 14
Set-TextCell "D14" "2.593.39"
Set-TextCell "E14" "  +0.22%  "

# Row 15
Set-TextCell "D15" "2.254.24"
Set-TextCell "E15" "  +0.26%  "

# Row 16
Set-TextCell "D16" "0.832"
Set-TextCell "E16" "  -2.59%  "

# Row 17
Set-TextCell "D17" "13.73"
Set-TextCell "E17" "  -2.78%  "

# Row 18
Set-TextCell "D18" "43.980.32"
Set-TextCell "E18" "  +0.82%  "

# Row 19
Set-TextCell "D19" "12.64"
Set-TextCell "E19" "  -6.96%  "

# Row 20
Set-TextCell "D20" "0.0₃0969"
Set-TextCell "E20" "  -1.68%  "

# Row 21
Set-TextCell "D21" "6.30"
Set-TextCell "E21" "  -3.75%  "

# Row 22
Set-TextCell "D22" "65.16"
Set-TextCell "E22" "  -0.21%  "

# Row 23
Set-TextCell "D23" "240.73"
Set-TextCell "E23" "  +1.81%  "

# Row 24
Set-TextCell "D24" "2.94"
Set-TextCell "E24" "  -7.83%  "

# Row 25
Set-TextCell "D25" "1.96"
Set-TextCell "E25" "  -9.14%  "

# Row 26
Set-TextCell "E26" "  +0.25%  "

# Row 27
Set-TextCell "D27" "10.05"
Set-TextCell "E27" "  -0.41%  "

# Row 28
Set-TextCell "E28" "  -2.16%  "

# Row 29
Set-TextCell "D29" "36.45"
Set-TextCell "E29" "  -1.15%  "

# Row 30
Set-TextCell "D30" "6.14"
Set-TextCell "E30" "  -2.35%  "

# Row 31
Set-TextCell "D31" "20.06"
Set-TextCell "E31" "  -0.46%  "

# Row 32
Set-TextCell "D32" "156.27"
Set-TextCell "E32" "  -2.74%  "

# Row 33
Set-TextCell "D33" "3.47"
Set-TextCell "E33" "  +13.73%  "

# Row 34
Set-TextCell "D34" "0.0817"
Set-TextCell "E34" "  -4.24%  "

# Row 35
Set-TextCell "E35" "  -0.98%  "

# Row 36
Set-TextCell "E36" "  +0.06%  "

# Row 37
Set-TextCell "D37" "0.106"
Set-TextCell "E37" "  -5.44%  "

# Row 38
Set-TextCell "D38" "1.82"
Set-TextCell "E38" "  -4.78%  "

# Row 39
Set-TextCell "D39" "15.36"
Set-TextCell "E39" "  -1.67%  "

# Row 40
Set-TextCell "B40" "NEARProtocol"
Set-TextCell "C40" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D40" "3.36"
Set-TextCell "E40" "  -10.78%  "

# Row 41
Set-TextCell "D41" "0.0305"
Set-TextCell "E41" "  -4.06%  "

# Row 42
Set-TextCell "B42" "RenderToken"
Set-TextCell "C42" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D42" "3.82"
Set-TextCell "E42" "  -9.82%  "

# Row 43
Set-TextCell "D43" "1.01"
Set-TextCell "E43" "  +0.15%  "

# Row 44
Set-TextCell "D44" "1.758.54"
Set-TextCell "E44" "  -2.08%  "

# Row 45
Set-TextCell "D45" "86.83"
Set-TextCell "E45" "  +5.49%  "

# Row 46
Set-TextCell "D46" "5.14"
Set-TextCell "E46" "  -1.10%  "

# Row 47
Set-TextCell "D47" "0.191"
Set-TextCell "E47" "  -3.76%  "

# Row 48
Set-TextCell "D48" "100.92"
Set-TextCell "E48" "  -2.36%  "

# Row 49
Set-TextCell "D49" "8.22"
Set-TextCell "E49" "  -2.42%  "

# Row 50
Set-TextCell "D50" "55.15"
Set-TextCell "E50" "  -5.99%  "

# Row 51
Set-TextCell "D51" "68.55"
Set-TextCell "E51" "  -9.66%  "
